# Fix minutes/seconds formatting in the "haul" (total time) column (I):
# zero-pad single-digit minutes and seconds, e.g. "11 ч. 22 мин. 9 сек."
# becomes "11 ч. 22 мин. 09 сек.". Hours are left unpadded.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count

for ($r = 1; $r -le $rowCount; $r++) {
  $cell = $ws.Cells.Item($r, 9)
  $val = $cell.Value()
  if ($val -ne $null) {
    $s = [string]$val
    if ($s -match "^(\d+) ч\. (\d+) мин\. (\d+) сек\.$") {
      $h = $matches[1]
      $mi = $matches[2]
      $se = $matches[3]
      $needsPad = ($mi.Length -lt 2) -or ($se.Length -lt 2)
      if ($needsPad) {
        $miP = $mi.PadLeft(2, "0")
        $seP = $se.PadLeft(2, "0")
        $new = "$h ч. $miP мин. $seP сек."
        $cell.Value = $new
      }
    }
  }
}
